# Updates cryptos list prices/volumes (and one row reorder: LidoDAOToken /
# EthereumClassic swap positions 28/29) per the source feed refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry gives the new "Price" (D) / "Volume(1h)" (E) text for a row,
# plus new "Coin" (B) / "Link" (C) text where a row's identity changed.
$updates = @(
    @{ Row=2; D="30.282.66"; E="  -3.41%  " }
    @{ Row=3; D="1.931.13"; E="  -3.26%  " }
    @{ Row=4; D="0.9990"; E="  -0.03%  " }
    @{ Row=5; D="247.15"; E="  -2.96%  " }
    @{ Row=6; D="0.7241"; E="  -10.46%  " }
    @{ Row=7; D="0.9986"; E="  -0.03%  " }
    @{ Row=8; D="0.3290"; E="  -6.65%  " }
    @{ Row=9; D="26.82"; E="  +4.87%  " }
    @{ Row=10; D="0.06834"; E="  -2.71%  " }
    @{ Row=11; D="0.8075"; E="  -4.04%  " }
    @{ Row=12; D="0.07962"; E="  -1.98%  " }
    @{ Row=13; D="1.930.99"; E="  -2.89%  " }
    @{ Row=14; D="5.437"; E="  -1.87%  " }
    @{ Row=15; D="94.76"; E="  -6.34%  " }
    @{ Row=16; D="14.59"; E="  +4.25%  " }
    @{ Row=17; D="262.19"; E="  -3.92%  " }
    @{ Row=18; D="30.275.03" }
    @{ Row=19; D="0.000007952"; E="  -0.10%  " }
    @{ Row=20; D="5.835"; E="  +0.44%  " }
    @{ Row=21; D="2.183.49"; E="  -3.04%  " }
    @{ Row=22; D="0.9990"; E="  +0.03%  " }
    @{ Row=23; D="0.9985"; E="  -0.06%  " }
    @{ Row=24; D="6.916"; E="  -0.80%  " }
    @{ Row=25; D="9.719"; E="  -1.07%  " }
    @{ Row=26; D="160.12"; E="  -2.41%  " }
    @{ Row=27; D="0.1358"; E="  -9.93%  " }
    @{ Row=28; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.333"; E="  +4.56%  " }
    @{ Row=29; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="19.03"; E="  -5.20%  " }
    @{ Row=30; D="1.364"; E="  +0.50%  " }
    @{ Row=31; D="1.551" }
    @{ Row=32; D="4.409"; E="  -3.97%  " }
    @{ Row=33; D="4.215"; E="  -2.92%  " }
    @{ Row=34; D="0.05100"; E="  -1.98%  " }
    @{ Row=35; D="1.208"; E="  -0.71%  " }
    @{ Row=36; D="0.7457"; E="  -1.84%  " }
    @{ Row=37; D="2.723"; E="  -1.83%  " }
    @{ Row=38; D="0.01947"; E="  -3.25%  " }
    @{ Row=39; D="2.818"; E="  -3.31%  " }
    @{ Row=40; D="80.45"; E="  +2.73%  " }
    @{ Row=41; D="6.596"; E="  -0.94%  " }
    @{ Row=42; D="0.4487"; E="  -5.79%  " }
    @{ Row=43; D="2.018"; E="  -4.19%  " }
    @{ Row=44; D="0.9991"; E="  +0.00%  " }
    @{ Row=45; D="0.8376"; E="  -2.22%  " }
    @{ Row=46; D="102.55"; E="  -1.78%  " }
    @{ Row=47; D="9.751"; E="  -2.25%  " }
    @{ Row=48; D="7.338"; E="  -2.50%  " }
    @{ Row=49; D="36.37"; E="  -1.57%  " }
    @{ Row=50; D="0.4134"; E="  -5.52%  " }
    @{ Row=51; D="1.488"; E="  +2.21%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $ws.Range("B$($u.Row)").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$($u.Row)").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Force text formatting so values like "0.9990" or "30.282.66" are
        # preserved exactly (matching the source sheet's inline-string cells)
        # instead of being auto-coerced into numbers and losing precision.
        $dCell = $ws.Range("D$($u.Row)")
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $eCell = $ws.Range("E$($u.Row)")
        $eCell.NumberFormat = "@"
        $eCell.Value = $u.E
    }
}
